$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("datasetUniform_dag10.xlsxbeginOrder", 0.8348260481567809, 0.165173951843219, 3.852659225463867),
    @("datasetUniform_dag10.xlsxendOrder", 0.09736735146161822, 0.9026326485383818, 3.835338115692139),
    @("datasetUniform_dag20.xlsxbeginOrder", 0.8239025968423785, 0.1760974031576215, 5.280167579650879),
    @("datasetUniform_dag20.xlsxendOrder", 0.4431607686553564, 0.5568392313446436, 8.93315577507019)
)

$startRow = 19
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
